$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added. It is inserted as row 6, pushing the
# previous rows 6-11 down to rows 7-12.
$ws.Rows.Item(6).Insert()

# Fill in the new row 6 with the new weekly data point.
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = "Terminal La Palmera de La Serena"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D6").Value = 44466
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100112026
$ws.Range("G6").Value = "Haba"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 9500
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 9750
$ws.Range("N6").Value = "`$/saco 25 kilos"
$ws.Range("O6").Value = "Provincia del Elquí"
$ws.Range("P6").Value = 390
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
